$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "cost minimization (Ref-A1B2C3)" "cost minimization (Ref-u315816)"
Replace-Text "absence of RL in strategic planning (Ref-AB1CD2)" "absence of RL in strategic planning (Ref-f559901)"
Replace-Text "Antonyová, et al. (Ref-A1B2C3), there are factors" "Antonyová, et al. (Ref-u733580), there are factors"
Replace-Text "Antonyová, et al. (Ref-D4E5F6), these trends" "Antonyová, et al. (Ref-f748538), these trends"
Replace-Text "Antonyová, et al. (Ref-J7X8A2). It considers" "Antonyová, et al. (Ref-s950403). It considers"
Replace-Text "economic growth of the country (Ref-AB1CD2). This model" "economic growth of the country (Ref-u136632). This model"
Replace-Text "repair location for resale (Ref-A1B2C3). Moreover" "repair location for resale (Ref-s457446). Moreover"
Replace-Text "can forecast data (Ref-SG7J2K)" "can forecast data (Ref-f511197)"
Replace-Text "Business Wire (Ref-K92Y6F)" "Business Wire (Ref-f511197)"
Replace-Text "amendments of federal solid waste (Ref-J7X2B9). It describe" "amendments of federal solid waste (Nguyen, 2015). It describe"
Replace-Text "guarantees no harm to the environment (Ref-A1B2C3). Therefo" "guarantees no harm to the environment (Ref-f970291). Therefo"
Replace-Text "impact of nicotine wastes (Ref-J7X2N9). Finally" "impact of nicotine wastes (Ref-u107756). Finally"
